$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17 (pushes existing rows 17-85 down to 18-86)
$ws.Rows(17).Insert()

# Populate the newly inserted row 17 with the new weekly record
$ws.Cells.Item(17, 1).Value = 9
$ws.Cells.Item(17, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(17, 3).Value = "Metropolitana"
$ws.Cells.Item(17, 4).Value = 45168
$ws.Cells.Item(17, 5).Value = 13
$ws.Cells.Item(17, 6).Value = 100112010
$ws.Cells.Item(17, 7).Value = "Achicoria"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 70
$ws.Cells.Item(17, 11).Value = 7000
$ws.Cells.Item(17, 12).Value = 7000
$ws.Cells.Item(17, 13).Value = 7000
$ws.Cells.Item(17, 14).Value = "`$/caja 16 unidades"
$ws.Cells.Item(17, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(17, 16).Value = 438
$ws.Cells.Item(17, 17).Value = 16
$ws.Cells.Item(17, 18).Value = "Hortaliza"
